$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.531.25'
$ws.Range('E2').Value = '  +3.74%  '
$ws.Range('D3').Value = '1.587.76'
$ws.Range('E3').Value = '  +1.02%  '
$ws.Range('E4').Value = '  +0.94%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.97'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.79%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.494'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.16%  '
$ws.Range('E7').Value = '  +0.92%  '
$ws.Range('E8').Value = '  +5.75%  '
$ws.Range('E9').Value = '  +0.60%  '
$ws.Range('E10').Value = '  +0.91%  '
$ws.Range('E11').Value = '  +1.66%  '
$ws.Range('D12').Value = '1.814.44'
$ws.Range('E12').Value = '  +1.01%  '
$ws.Range('D13').Value = '1.586.95'
$ws.Range('E13').Value = '  +1.11%  '
$ws.Range('E14').Value = '  +1.66%  '
$ws.Range('E15').Value = '  -0.48%  '
$ws.Range('D16').Value = '28.530.44'
$ws.Range('E16').Value = '  +3.87%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.08'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '231.04'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.14%  '
$ws.Range('E19').Value = '  -0.22%  '
$ws.Range('D20').Value = '0.0₃0706'
$ws.Range('E20').Value = '  +0.19%  '
$ws.Range('E21').Value = '  +0.93%  '
$ws.Range('E22').Value = '  -1.78%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.33'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.82%  '
$ws.Range('E24').Value = '  +2.31%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.80'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.86%  '
$ws.Range('E26').Value = '  +0.30%  '
$ws.Range('E27').Value = '  -0.88%  '
$ws.Range('E28').Value = '  -0.60%  '
$ws.Range('E29').Value = '  +0.87%  '
$ws.Range('E30').Value = '  -0.98%  '
$ws.Range('E31').Value = '  -0.60%  '
$ws.Range('E32').Value = '  +0.20%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.18'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.36%  '
$ws.Range('D34').Value = '1.392.32'
$ws.Range('E34').Value = '  -4.39%  '
$ws.Range('E35').Value = '  -1.34%  '
$ws.Range('E36').Value = '  -10.58%  '
$ws.Range('E37').Value = '  +1.07%  '
$ws.Range('E38').Value = '  +10.59%  '
$ws.Range('E39').Value = '  -0.76%  '
$ws.Range('E40').Value = '  -0.15%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.812'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.18%  '
$ws.Range('E42').Value = '  +0.86%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.65'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.13%  '
$ws.Range('E44').Value = '  +0.86%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.982'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.02%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '62.95'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.55%  '
$ws.Range('D47').Value = '1.723.99'
$ws.Range('E47').Value = '  +0.94%  '
$ws.Range('E48').Value = '  +0.95%  '
$ws.Range('E49').Value = '  +0.02%  '
$ws.Range('D50').Value = '0.0₆0104'
$ws.Range('E50').Value = '  -0.80%  '
$ws.Range('E51').Value = '  -1.09%  '
